$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("General Information")
$ws2 = $wb.Worksheets.Item("Duration surgeries")

# Split the existing department groupings ("b" and "c") into finer
# sub-groups ("b"/"c" and "e"/"f"/"g") so the output is more granular
# (per the commit message: "Apart soort expanded output, print voor
# department").
$ws2.Range("C7:C11").Value = "b"
$ws2.Range("C12:C16").Value = "c"
$ws2.Range("C17:C21").Value = "e"
$ws2.Range("C22:C26").Value = "f"
$ws2.Range("C27:C31").Value = "g"

# Update the view/selection state on both sheets, activating the
# "Duration surgeries" sheet first so its selection can be recorded,
# then returning to "General Information" last so it remains the
# selected tab (as in the original workbook).
$ws2.Activate()
$ws2.Range("C31").Select()

$ws1.Activate()
$ws1.Range("F4").Select()
